$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.992.80"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").Value = "2.244.82"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'232.72"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").Value = "'0.636"
$ws.Range("E6").Value = "  -5.93%  "
$ws.Range("D7").Value = "'69.87"
$ws.Range("E7").Value = "  -4.80%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.559"
$ws.Range("E9").Value = "  -7.40%  "
$ws.Range("D10").Value = "'0.0991"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "'58.39"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "'35.77"
$ws.Range("E12").Value = "  +6.13%  "
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "'6.82"
$ws.Range("E14").Value = "  -6.83%  "
$ws.Range("D15").Value = "2.574.92"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").Value = "'15.14"
$ws.Range("E16").Value = "  -7.92%  "
$ws.Range("D17").Value = "'0.864"
$ws.Range("E17").Value = "  -5.12%  "
$ws.Range("D18").Value = "2.243.41"
$ws.Range("E18").Value = "  -4.81%  "
$ws.Range("D19").Value = "41.879.91"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  -4.57%  "
$ws.Range("D21").Value = "'6.26"
$ws.Range("E21").Value = "  -6.79%  "
$ws.Range("D22").Value = "'73.53"
$ws.Range("E22").Value = "  -5.36%  "
$ws.Range("D23").Value = "'237.88"
$ws.Range("E23").Value = "  -7.45%  "
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +5.00%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("D30").Value = "'169.38"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'20.68"
$ws.Range("E31").Value = "  -8.96%  "
$ws.Range("E32").Value = "  -7.19%  "
$ws.Range("E33").Value = "  -7.06%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("D36").Value = "'4.79"
$ws.Range("E36").Value = "  -8.19%  "
$ws.Range("D37").Value = "'3.61"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").Value = "'22.11"
$ws.Range("E38").Value = "  +16.00%  "
$ws.Range("D39").Value = "'2.27"
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("D40").Value = "'6.05"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").Value = "'67.22"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'9.10"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'4.91"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").Value = "'0.101"
$ws.Range("E45").Value = "  -9.67%  "
$ws.Range("D46").Value = "'0.190"
$ws.Range("E46").Value = "  -6.14%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  -5.24%  "
$ws.Range("D49").Value = "'4.37"
$ws.Range("E49").Value = "  +6.63%  "
$ws.Range("E50").Value = "  -7.01%  "
$ws.Range("D51").Value = "'9.94"
$ws.Range("E51").Value = "  +2.40%  "

# Reset style on forced-text numeric cells so no stray cell style/format is left behind
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
